$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adults")

# --- New column A width (diff: <col min="1" max="1" width="12.83203125" customWidth="1"/>) ---
$ws.Columns.Item(1).ColumnWidth = 12

# --- Row 14 height (diff: ht="256") ---
$ws.Rows.Item(14).RowHeight = 256

# --- Row 14 data ---
# The order below controls the order new strings are appended to sharedStrings.xml,
# matching the target document (A, D, E, F, Y, X, U, V, W).
$ws.Range("A14").Value = "Polio "
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "Conditional"
$ws.Range("D14").Value = "The first dose at any time"
$ws.Range("E14").Value = "1 or 2 months Spacing"
$ws.Range("F14").Value = "6 to 12 months Spacing"
$ws.Range("G14").Value = "X"
$ws.Range("H14").Value = "X"

$ws.Range("I14").Value = 6935
$ws.Range("J14").Value = 364635
$ws.Range("K14").Value = 6935
$ws.Range("L14").Value = 364635
$ws.Range("M14").Value = 6935
$ws.Range("N14").Value = 364635
$ws.Range("O14").Value = 6935
$ws.Range("P14").Value = 364635

$y14 = $ws.Range("Y14")
$y14.Value = "administer remaining doses"

$x14 = $ws.Range("X14")
$x14.WrapText = $true
$x14.Value = "No evidence of a complete polio vaccination series (i.e., at least 3 doses)`n(1, 2, or 3 doses) to complete a 3-dose series"

$u14 = $ws.Range("U14")
$u14.WrapText = $true
$u14.Value = "Routine poliovirus vaccination of adults residing in the United States is not necessary. It is recommended to adults at increased risk of exposure`nto poliovirus."

$v14 = $ws.Range("V14")
$v14.WrapText = $true
$v14.Value = "'adults at increased risk of exposure`nto poliovirus that have vidence of completed polio vaccination series`n(i.e., at least 3 doses)"

$w14 = $ws.Range("W14")
$w14.Value = "may administer one lifetime IPV booster"

# --- Selection (diff: activeCell="V18" sqref="V18") ---
$ws.Range("V18").Select()
